$d = $word.ActiveDocument

# 1) "The model consist on 6 subsystems..." -> "The model consists of 6 subsystems..."
$r1 = $d.Content
$r1.Find.ClearFormatting()
$found1 = $r1.Find.Execute("The model consist on", $false, $false, $false, $false, $false, $true, 1, $false, "The model consists of", 2)
Write-Host "Replaced 'consist on' -> 'consists of': $found1"

# 2) "This blocks contains the input..." -> "This block contains the input..."
$r2 = $d.Content
$r2.Find.ClearFormatting()
$found2 = $r2.Find.Execute("This blocks contains", $false, $false, $false, $false, $false, $true, 1, $false, "This block contains", 2)
Write-Host "Replaced 'This blocks contains' -> 'This block contains': $found2"

# 3) Move the "_GoBack" bookmark (Word's last-edit-location marker) from the
#    title line to right after "This block", reflecting that this was the
#    last text edited in the document.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$locate = $d.Content
$locate.Find.ClearFormatting()
$locate.Find.Execute("This block", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$newSpot = $d.Range($locate.End, $locate.End)
$d.Bookmarks.Add("_GoBack", $newSpot)
Write-Host "Moved _GoBack bookmark to $($newSpot.Start)"
